$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 184 (shifts existing rows 184:278 down to 185:279)
$ws.Rows(184).Insert()

# Populate the newly inserted row 184 with the new record's data
$ws.Range("A184").Value = 10
$ws.Range("B184").Value = "Vega Modelo de Temuco"
$ws.Range("C184").Value = "La Araucanía"
$ws.Range("D184").Value = 44523
$ws.Range("E184").Value = 9
$ws.Range("F184").Value = 100112024
$ws.Range("G184").Value = "Choclo"
$ws.Range("H184").Value = "Dulce o Americano"
$ws.Range("I184").Value = "Primera"
$ws.Range("J184").Value = 500
$ws.Range("K184").Value = 500
$ws.Range("L184").Value = 500
$ws.Range("M184").Value = 500
$ws.Range("N184").Value = '$/unidad'
$ws.Range("O184").Value = "Argentina"
$ws.Range("P184").Value = 500
$ws.Range("Q184").Value = 1
$ws.Range("R184").Value = "Hortaliza"

# Ensure date cell keeps the expected number format / style (s="2")
$ws.Range("D184").NumberFormat = $ws.Range("D185").NumberFormat
